$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the floating point value in A2 (slight recalculation drift)
$ws.Cells.Item(2, 1).Value = 45877.04185537037

# Add new row 3 with the latest weather reading
$ws.Cells.Item(3, 1).Value = 45877.08352210094
$ws.Cells.Item(3, 1).NumberFormat = $ws.Cells.Item(2, 1).NumberFormat
$ws.Cells.Item(3, 2).Value = 2025
$ws.Cells.Item(3, 3).Value = 32
$ws.Cells.Item(3, 4).Value = 13.78
$ws.Cells.Item(3, 5).Value = 92
$ws.Cells.Item(3, 6).Value = 0
$ws.Cells.Item(3, 7).Value = 7.23
$ws.Cells.Item(3, 8).Value = "SE"
$ws.Cells.Item(3, 9).Value = 0
$ws.Cells.Item(3, 10).Value = "02:00:16"
